$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the header style (bold font, thin border, centered alignment) from an
# existing header cell (AC1) onto the three new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels for team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player data row (2-55)
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 30).Value = 62   # AD = column 30 -> Wins
    $ws.Cells.Item($r, 31).Value = 100  # AE = column 31 -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32 -> Ties
}
